$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 10 (Objetivos PT) text: was the professor name placeholder, now the actual PT objectives paragraph ---
$ws.Range("B10").Value = "Como parte fundamental da formação específica e geral, a disciplina tem por objetivos (a) fornecer os conceitos fundamentais sobre dispositivos semicondutores aplicados em circuitos eletrônicos, que são importantes para a formação em engenharia física; (b) capacitar o aluno, trabalhando individualmente e em grupo, a modelar e resolver problemas de interesse envolvendo os principais componentes eletrônicos, como diodos, transistores e amplificadores operacionais, com escolhas adequadas de hipóteses e aplicação de ferramentas correspondentes de solução; (c) introduzir os componentes, técnicas, softwares e equipamentos utilizados na análise e projeto de circuitos eletrônicos; e (d) aplicar e estender os conceitos físicos aprendidos previamente."
$ws.Range("C10").Value = "Como parte fundamental da formação específica e geral, a disciplina tem por objetivos (a) fornecer os conceitos fundamentais sobre dispositivos semicondutores aplicados em circuitos eletrônicos, que são importantes para a formação em engenharia física; (b) capacitar o aluno, trabalhando individualmente e em grupo, a modelar e resolver problemas de interesse envolvendo os principais componentes eletrônicos, como diodos, transistores e amplificadores operacionais, com escolhas adequadas de hipóteses e aplicação de ferramentas correspondentes de solução; (c) introduzir os componentes, técnicas, softwares e equipamentos utilizados na análise e projeto de circuitos eletrônicos; e (d) aplicar e estender os conceitos físicos aprendidos previamente."

# --- Clear rows 13-23 entirely (content + formatting) so stray cells do not linger when the table is rebuilt ---
$ws.Range("A13:C23").Clear()

# --- Re-populate rows 13-25 per the corrected / expanded table layout ---
# Row 13
$ws.Range("B13").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C13").Value = "519033 - Carlos Yujiro Shigue"

# Row 14
$ws.Range("B14").Value = "7290967 - Emerson Gonçalves de Melo"
$ws.Range("C14").Value = "7290967 - Emerson Gonçalves de Melo"

# Row 15
$ws.Range("A15").Value = "Programa resumido:"
$ws.Range("B15").Value = "Introdução aos materiais e dispositivos semicondutores; Diodo; Transistor bipolar de junção; Transistores de efeito de campo; Amplificadores operacionais; Amplificadores de pequenos sinais; Fontes de alimentação. Análise e projeto de circuitos eletrônicos utilizando softwares EDA. Análises experimentais de circuitos eletrônicos."
$ws.Range("C15").Value = "Introdução aos materiais e dispositivos semicondutores; Diodo; Transistor bipolar de junção; Transistores de efeito de campo; Amplificadores operacionais; Amplificadores de pequenos sinais; Fontes de alimentação. Análise e projeto de circuitos eletrônicos utilizando softwares EDA. Análises experimentais de circuitos eletrônicos."

# Row 16
$ws.Range("A16").Value = "Short syllabus:"
$ws.Range("B16").Value = "Introduction to semiconductor materials and devices; Diode; Bipolar junction transistor; Field effect transistors; operational amplifiers; Small signal amplifiers; Power supplies. Analysis and design of electronic circuits using EDA software. Experimental analysis of electronic circuits."
$ws.Range("C16").Value = "Introduction to semiconductor materials and devices; Diode; Bipolar junction transistor; Field effect transistors; operational amplifiers; Small signal amplifiers; Power supplies. Analysis and design of electronic circuits using EDA software. Experimental analysis of electronic circuits."

# Row 17
$ws.Range("A17").Value = "Programa:"
$ws.Range("B17").Value = "1. Materiais Semicondutores; Diodos;    2. Retificadores de Tensão: Análise e Projeto;    3. Software EDA; Projeto de Placas de Circuito Impresso.    4. Transistor Bipolar de Junção (TBJ); Folha de Dados, Polarização e Chaveamento.    5. Transistor de Efeito de Campo (JFET - MOSFET); Folha de Dados, Polarização e Chaveamento.    6. Amplificadores para Pequenos Sinais: Análise e Projeto;    7. Amplificadores Operacionais: Buffer, Amplificação, Integração e Diferenciação;    8. Filtros Ativos: Análise e Projeto;    9. Fontes de Alimentação;    10. Circuitos Optoeletrônicos;"
$ws.Range("C17").Value = "1. Materiais Semicondutores; Diodos;    2. Retificadores de Tensão: Análise e Projeto;    3. Software EDA; Projeto de Placas de Circuito Impresso.    4. Transistor Bipolar de Junção (TBJ); Folha de Dados, Polarização e Chaveamento.    5. Transistor de Efeito de Campo (JFET - MOSFET); Folha de Dados, Polarização e Chaveamento.    6. Amplificadores para Pequenos Sinais: Análise e Projeto;    7. Amplificadores Operacionais: Buffer, Amplificação, Integração e Diferenciação;    8. Filtros Ativos: Análise e Projeto;    9. Fontes de Alimentação;    10. Circuitos Optoeletrônicos;"

# Row 18
$ws.Range("A18").Value = "Syllabus:"
$ws.Range("B18").Value = "1. Semiconductor Materials; Diodes;2. Voltage Rectifiers: Analysis and Design;3. EDA Software; Design of Printed Circuit Boards.4. Bipolar Junction Transistor (BJT); Data Sheet, Polarization and Switching.5. Field Effect Transistor (JFET - MOSFET); Data Sheet, Polarization and Switching.6. Small Signal Amplifiers: Analysis and Design;7. Operational Amplifiers: Buffer, Amplification, Integration and Differentiation;8. Active Filters: Analysis and Design;9. Power Supplies;10. Optoelectronic Circuits;"
$ws.Range("C18").Value = "1. Semiconductor Materials; Diodes;2. Voltage Rectifiers: Analysis and Design;3. EDA Software; Design of Printed Circuit Boards.4. Bipolar Junction Transistor (BJT); Data Sheet, Polarization and Switching.5. Field Effect Transistor (JFET - MOSFET); Data Sheet, Polarization and Switching.6. Small Signal Amplifiers: Analysis and Design;7. Operational Amplifiers: Buffer, Amplification, Integration and Differentiation;8. Active Filters: Analysis and Design;9. Power Supplies;10. Optoelectronic Circuits;"

# Row 19
$ws.Range("A19").Value = "Avaliação:"

# Row 20
$ws.Range("A20").Value = "Método:"
$ws.Range("B20").Value = "Aulas expositivas e práticas de laboratório com interações em grupo para a solução de problemas."
$ws.Range("C20").Value = "Aulas expositivas e práticas de laboratório com interações em grupo para a solução de problemas."

# Row 21
$ws.Range("A21").Value = "Critério:"
$ws.Range("B21").Value = "Média aritmética (M) de provas individuais (P1 e P2) e trabalhos em grupo ao longo do semestre (T), tal que M = 0,3*P1+0,3*P2+0.4*T"
$ws.Range("C21").Value = "Média aritmética (M) de provas individuais (P1 e P2) e trabalhos em grupo ao longo do semestre (T), tal que M = 0,3*P1+0,3*P2+0.4*T"

# Row 22
$ws.Range("A22").Value = "Norma de recuperação:"
$ws.Range("B22").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Range("C22").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"

# Row 23
$ws.Range("A23").Value = "Bibliografia:"
$ws.Range("B23").Value = "BROPHY, J. J. Eletrônica Básica. Guanabara Dois.NOVO, D. D. Eletrônica Aplicada. Editora da USP.SIMPSON, R.E. Introductory electronics for scientists and engineers. Allyn and Bacon.HOROWITZ, P.; HILL, W. The art of electronics. Cambridge University Press.MOTCHENBACHER, C. D.; FITCHEN, F.C. Low noise electronic design, John Wiley and Sons.MORRISON, R. Grounding and shielding techniques in instrumentation, John Wiley and Sons.ALEXANDER, C. K. E SADIKU, M. N. O. Fundamentos de Circuitos Elétricos. McGraw-Hill, 2013.NILSSON, J. W. E RIEDEL, S. A. Electric Circuits. Prentice Hall, 2011.BOYLESTAD, R. L. E NASHELSKY, L. Electronic Devices andCircuit Theory. Pearson, 2013"
$ws.Range("C23").Value = "BROPHY, J. J. Eletrônica Básica. Guanabara Dois.NOVO, D. D. Eletrônica Aplicada. Editora da USP.SIMPSON, R.E. Introductory electronics for scientists and engineers. Allyn and Bacon.HOROWITZ, P.; HILL, W. The art of electronics. Cambridge University Press.MOTCHENBACHER, C. D.; FITCHEN, F.C. Low noise electronic design, John Wiley and Sons.MORRISON, R. Grounding and shielding techniques in instrumentation, John Wiley and Sons.ALEXANDER, C. K. E SADIKU, M. N. O. Fundamentos de Circuitos Elétricos. McGraw-Hill, 2013.NILSSON, J. W. E RIEDEL, S. A. Electric Circuits. Prentice Hall, 2011.BOYLESTAD, R. L. E NASHELSKY, L. Electronic Devices andCircuit Theory. Pearson, 2013"

# Row 24
$ws.Range("A24").Value = "Requisitos:"

# Row 25
$ws.Range("B25").Value = "LOM3262 -  Circuitos Elétricos  (Requisito)`n"
$ws.Range("C25").Value = "LOM3262 -  Circuitos Elétricos  (Requisito)`n"

# --- Apply the same per-column formatting (font/wrap) used throughout the sheet to the newly written cells ---
$ws.Range("A3").Copy()
$ws.Range("A13:A25").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B13:B25").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C13:C25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row heights ---
$ws.Rows.Item(13).RowHeight = 15
$ws.Rows.Item(14).RowHeight = 15
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 60
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 120
$ws.Rows.Item(19).RowHeight = 15
$ws.Rows.Item(20).RowHeight = 15
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 60
$ws.Rows.Item(23).RowHeight = 120
$ws.Rows.Item(24).RowHeight = 15
$ws.Rows.Item(25).RowHeight = 30

# --- Clear any stray A/B/C cells left over on rows where that column should now be empty ---
$ws.Range("A13").Clear()
$ws.Range("A14").Clear()
$ws.Range("B19").Clear()
$ws.Range("C19").Clear()
$ws.Range("B24").Clear()
$ws.Range("C24").Clear()
$ws.Range("A25").Clear()

# --- Column A used to be merged with column B in one <col min="1" max="2"> span; split it so col A is its own span (col B keeps its own width/style) ---
$ws.Columns.Item(1).Hidden = $false
